$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-19 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-20 Sunday", 2)
$d.Content.Find.Execute("93×95=8835", $true, $false, $false, $false, $false, $true, 1, $false, "92×77=7084", 2)
$d.Content.Find.Execute("45×95=4275", $true, $false, $false, $false, $false, $true, 1, $false, "49×19=931", 2)
$d.Content.Find.Execute("13×99=1287", $true, $false, $false, $false, $false, $true, 1, $false, "69×45=3105", 2)
$d.Content.Find.Execute("50×24=1200", $true, $false, $false, $false, $false, $true, 1, $false, "72×28=2016", 2)
$d.Content.Find.Execute("26×69=1794", $true, $false, $false, $false, $false, $true, 1, $false, "18×97=1746", 2)
$d.Content.Find.Execute("46×27=1242", $true, $false, $false, $false, $false, $true, 1, $false, "36×72=2592", 2)
$d.Content.Find.Execute("65×55=3575", $true, $false, $false, $false, $false, $true, 1, $false, "59×57=3363", 2)
$d.Content.Find.Execute("75×41=3075", $true, $false, $false, $false, $false, $true, 1, $false, "35×15=525", 2)
$d.Content.Find.Execute("63×24=1512", $true, $false, $false, $false, $false, $true, 1, $false, "28×41=1148", 2)
$d.Content.Find.Execute("64×69=4416", $true, $false, $false, $false, $false, $true, 1, $false, "86×14=1204", 2)
$d.Content.Find.Execute("52×18=936", $true, $false, $false, $false, $false, $true, 1, $false, "56×16=896", 2)
$d.Content.Find.Execute("76×83=6308", $true, $false, $false, $false, $false, $true, 1, $false, "82×71=5822", 2)
$d.Content.Find.Execute("71×47=3337", $true, $false, $false, $false, $false, $true, 1, $false, "65×79=5135", 2)
$d.Content.Find.Execute("72×73=5256", $true, $false, $false, $false, $false, $true, 1, $false, "97×37=3589", 2)
$d.Content.Find.Execute("65×23=1495", $true, $false, $false, $false, $false, $true, 1, $false, "39×15=585", 2)
$d.Content.Find.Execute("37×20=740", $true, $false, $false, $false, $false, $true, 1, $false, "33×48=1584", 2)
$d.Content.Find.Execute("57×14=798", $true, $false, $false, $false, $false, $true, 1, $false, "49×30=1470", 2)
$d.Content.Find.Execute("66×88=5808", $true, $false, $false, $false, $false, $true, 1, $false, "48×93=4464", 2)
$d.Content.Find.Execute("38×48=1824", $true, $false, $false, $false, $false, $true, 1, $false, "64×90=5760", 2)
$d.Content.Find.Execute("32×73=2336", $true, $false, $false, $false, $false, $true, 1, $false, "61×84=5124", 2)
$d.Content.Find.Execute("85×73=6205", $true, $false, $false, $false, $false, $true, 1, $false, "50×14=700", 2)
$d.Content.Find.Execute("36×79=2844", $true, $false, $false, $false, $false, $true, 1, $false, "77×72=5544", 2)
$d.Content.Find.Execute("34×11=374", $true, $false, $false, $false, $false, $true, 1, $false, "63×68=4284", 2)
$d.Content.Find.Execute("97×40=3880", $true, $false, $false, $false, $false, $true, 1, $false, "22×55=1210", 2)
$d.Content.Find.Execute("40×90=3600", $true, $false, $false, $false, $false, $true, 1, $false, "74×65=4810", 2)
